# Populate the "strain" column (F) with strain names that correspond to
# the existing "genotype" column (G) values, per row-group:
#   CNAG_00000                -> KN99alpha   (rows 2-4, 14-16)
#   CNAG_07797                -> TDY1256     (rows 5-7)
#   CNAG_01523.CNAG_01551     -> TDY1373     (rows 8-10)
#   CNAG_01551.CNAG_02153     -> TDY1367     (rows 11-13)
#   CNAG_05431                -> TDY1208     (rows 17-19)
#   CNAG_00440                -> TDY1210     (rows 20-22)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the first occurrence of each new strain name in the same order
# they were first introduced in the authored workbook, so the shared
# string table is built up in that order: TDY1373, TDY1367, KN99alpha,
# TDY1256, TDY1208, TDY1210.
$ws.Range("F8").Value = "TDY1373"
$ws.Range("F11").Value = "TDY1367"
$ws.Range("F2").Value = "KN99alpha"
$ws.Range("F5").Value = "TDY1256"
$ws.Range("F17").Value = "TDY1208"
$ws.Range("F20").Value = "TDY1210"

# Fill in the remaining rows of each group.
$ws.Range("F3").Value = "KN99alpha"
$ws.Range("F4").Value = "KN99alpha"

$ws.Range("F6").Value = "TDY1256"
$ws.Range("F7").Value = "TDY1256"

$ws.Range("F9").Value = "TDY1373"
$ws.Range("F10").Value = "TDY1373"

$ws.Range("F12").Value = "TDY1367"
$ws.Range("F13").Value = "TDY1367"

$ws.Range("F14").Value = "KN99alpha"
$ws.Range("F15").Value = "KN99alpha"
$ws.Range("F16").Value = "KN99alpha"

$ws.Range("F18").Value = "TDY1208"
$ws.Range("F19").Value = "TDY1208"

$ws.Range("F21").Value = "TDY1210"
$ws.Range("F22").Value = "TDY1210"

# Match the author's final on-screen state: scrolled so row 11 is at the
# top, with F21:F22 selected (active cell F21).
[void]$ws.Range("F21:F22").Select()
